$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited columns (D = Price, E = Volume, plus B/C for the two swapped rows)
# hold text-formatted values in the source sheet (e.g. "1.00", "0.0000156" must
# keep trailing zeros / avoid scientific notation), so force text format before
# assigning -- mirrors how the sheet was originally authored (inline text cells).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.639.23'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.606.93'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.31'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.44'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.73%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.110'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.69'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.386'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.42%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.66'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.076.25'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.501.62'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000156'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +7.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.602.78'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.53'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +8.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.75'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '347.51'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.95'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.51'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +6.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.35'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.70'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '573.27'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.07'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.42%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.06'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0856'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.77'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.26'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '167.22'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.415'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.32%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.67'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.95'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '168.51'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.66'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.97'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +5.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0590'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.21'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.633'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.06'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.89%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0253'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0966'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.22'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0237'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +21.04%  '
